$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# -------------------------------------------------------------------
# Settings sheet (sheet1.xml)
# -------------------------------------------------------------------

# Row 2: OrchestratorQueueName value changes from "ProcessABCQueue" to "ACME_Client"
$wsSettings.Cells.Item(2, 2).Value = "ACME_Client"

# Row 3: OrchestratorAssetFolder value gets populated
$wsSettings.Cells.Item(3, 2).Value = "REF_Generate Yearly Report_Performer"

# Row 9: Status / Completed (B9 gets left-aligned)
$wsSettings.Cells.Item(9, 1).Value = "Status"
$wsSettings.Cells.Item(9, 2).Value = "Completed"
$wsSettings.Cells.Item(9, 2).HorizontalAlignment = -4131

# Row 10: Year_Report / 2023 (B10 gets left-aligned)
$wsSettings.Cells.Item(10, 1).Value = "Year_Report"
$wsSettings.Cells.Item(10, 2).Value = 2023
$wsSettings.Cells.Item(10, 2).HorizontalAlignment = -4131

# Row 11: URL section header (large bold font)
$wsSettings.Cells.Item(11, 1).Value = "URL"
$wsSettings.Cells.Item(11, 1).Font.Bold = $true
$wsSettings.Cells.Item(11, 1).Font.Size = 18
$wsSettings.Cells.Item(11, 1).Font.Name = "Calibri"

# Row 12: ACME_URL
$wsSettings.Cells.Item(12, 1).Value = "ACME_URL"
$wsSettings.Cells.Item(12, 2).Value = "https://acme-test.uipath.com/"

# Row 13: ACME_Workitem_URL
$wsSettings.Cells.Item(13, 1).Value = "ACME_Workitem_URL"
$wsSettings.Cells.Item(13, 2).Value = "https://acme-test.uipath.com/work-items/"

# Row 14: Download Monthly Report_URL
$wsSettings.Cells.Item(14, 1).Value = "Download Monthly Report_URL"
$wsSettings.Cells.Item(14, 2).Value = "https://acme-test.uipath.com/reports/download/"

# Row 15: Upload Yearly Report_URL
$wsSettings.Cells.Item(15, 1).Value = "Upload Yearly Report_URL"
$wsSettings.Cells.Item(15, 2).Value = "https://acme-test.uipath.com/reports/upload/"

# Row 18: Path section header (large bold font)
$wsSettings.Cells.Item(18, 1).Value = "Path"
$wsSettings.Cells.Item(18, 1).Font.Bold = $true
$wsSettings.Cells.Item(18, 1).Font.Size = 18
$wsSettings.Cells.Item(18, 1).Font.Name = "Calibri"
$wsSettings.Cells.Item(18, 1).Font.Color = 0

# Row 19: Report_Data_Download
$wsSettings.Cells.Item(19, 1).Value = "Report_Data_Download"
$wsSettings.Cells.Item(19, 2).Value = "Data\Report\"

# Remove the now-unused trailing empty rows 990:998 (dimension shrinks to Z989)
$wsSettings.Range("A990:A998").EntireRow.Delete()

# -------------------------------------------------------------------
# Constants sheet (sheet2.xml)
# -------------------------------------------------------------------

# Row 2: MaxRetryNumber value changes from 0 to 2
$wsConstants.Cells.Item(2, 2).Value = 2

# -------------------------------------------------------------------
# Selections / active sheet
# -------------------------------------------------------------------

# Set selection on Assets and Constants first (selecting a range also
# activates that sheet), then finish on Settings so it ends up active.
$wsAssets.Range("A1").Select()

$wsConstants.Range("C24").Select()

$wsSettings.Activate()
$wsSettings.Range("C20:C21").Select()
